# Insert a new weekly price record as row 265 (shifts existing rows 265-337 down to 266-338)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(265).Insert()

$newRow = 265

$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 45093
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100108
$ws.Cells.Item($newRow, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value = 100108005
$ws.Cells.Item($newRow, 10).Value = "Piña"
$ws.Cells.Item($newRow, 11).Value = "Caramelo"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 60
$ws.Cells.Item($newRow, 14).Value = 22000
$ws.Cells.Item($newRow, 15).Value = 22000
$ws.Cells.Item($newRow, 16).Value = 22000
$ws.Cells.Item($newRow, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item($newRow, 18).Value = "Ecuador"
$ws.Cells.Item($newRow, 19).Value = 1833
$ws.Cells.Item($newRow, 20).Value = 12
